{"js": "// Replace each three-digit-divided-by-one-digit division problem text\n// with its updated counterpart, per the authoritative old->new mapping.\n// All old values are unique text runs within the document body, so an\n// exact-text search/replace per pair is safe and order independent.\nconst replacements = [\n  [\"674\u00f73=\", \"704\u00f78=\"],\n  [\"218\u00f79=\", \"523\u00f77=\"],\n  [\"188\u00f77=\", \"935\u00f74=\"],\n  [\"426\u00f74=\", \"391\u00f77=\"],\n  [\"624\u00f77=\", \"439\u00f73=\"],\n  [\"312\u00f79=\", \"943\u00f73=\"],\n  [\"116\u00f74=\", \"743\u00f78=\"],\n  [\"438\u00f73=\", \"533\u00f79=\"],\n  [\"595\u00f73=\", \"200\u00f78=\"],\n  [\"644\u00f75=\", \"607\u00f75=\"],\n  [\"720\u00f74=\", \"234\u00f76=\"],\n  [\"569\u00f76=\", \"134\u00f79=\"],\n  [\"313\u00f75=\", \"476\u00f73=\"],\n  [\"986\u00f75=\", \"178\u00f73=\"],\n  [\"714\u00f79=\", \"579\u00f77=\"],\n  [\"874\u00f72=\", \"723\u00f79=\"],\n  [\"945\u00f75=\", \"925\u00f74=\"],\n  [\"637\u00f79=\", \"280\u00f77=\"],\n  [\"793\u00f73=\", \"165\u00f73=\"],\n  [\"156\u00f76=\", \"117\u00f73=\"],\n  [\"373\u00f78=\", \"196\u00f78=\"],\n  [\"110\u00f72=\", \"886\u00f73=\"],\n  [\"159\u00f74=\", \"235\u00f75=\"],\n  [\"764\u00f77=\", \"911\u00f75=\"],\n  [\"881\u00f75=\", \"485\u00f74=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old -> new text for each division problem in the practice table.\n# Every old value is a unique text run in the document, so a plain\n# Find/Replace-All pass per pair is safe and order independent.\n$replacements = @(\n    @(\"674\u00f73=\", \"704\u00f78=\"),\n    @(\"218\u00f79=\", \"523\u00f77=\"),\n    @(\"188\u00f77=\", \"935\u00f74=\"),\n    @(\"426\u00f74=\", \"391\u00f77=\"),\n    @(\"624\u00f77=\", \"439\u00f73=\"),\n    @(\"312\u00f79=\", \"943\u00f73=\"),\n    @(\"116\u00f74=\", \"743\u00f78=\"),\n    @(\"438\u00f73=\", \"533\u00f79=\"),\n    @(\"595\u00f73=\", \"200\u00f78=\"),\n    @(\"644\u00f75=\", \"607\u00f75=\"),\n    @(\"720\u00f74=\", \"234\u00f76=\"),\n    @(\"569\u00f76=\", \"134\u00f79=\"),\n    @(\"313\u00f75=\", \"476\u00f73=\"),\n    @(\"986\u00f75=\", \"178\u00f73=\"),\n    @(\"714\u00f79=\", \"579\u00f77=\"),\n    @(\"874\u00f72=\", \"723\u00f79=\"),\n    @(\"945\u00f75=\", \"925\u00f74=\"),\n    @(\"637\u00f79=\", \"280\u00f77=\"),\n    @(\"793\u00f73=\", \"165\u00f73=\"),\n    @(\"156\u00f76=\", \"117\u00f73=\"),\n    @(\"373\u00f78=\", \"196\u00f78=\"),\n    @(\"110\u00f72=\", \"886\u00f73=\"),\n    @(\"159\u00f74=\", \"235\u00f75=\"),\n    @(\"764\u00f77=\", \"911\u00f75=\"),\n    @(\"881\u00f75=\", \"485\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
